$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '26.960.70'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = "'" + '1.818.17'
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = "'" + '309.83'
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").Value = "'" + '0.4656'
$ws.Range("E7").Value = '  +0.82%  '
$ws.Range("D8").Value = "'" + '0.3662'
$ws.Range("E8").Value = '  -1.35%  '
$ws.Range("D9").Value = "'" + '0.07367'
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = "'" + '0.8718'
$ws.Range("E10").Value = '  -0.46%  '
$ws.Range("E11").Value = '  -1.07%  '
$ws.Range("D12").Value = "'" + '1.824.22'
$ws.Range("E12").Value = '  +2.55%  '
$ws.Range("D13").Value = "'" + '5.390'
$ws.Range("E13").Value = '  +0.51%  '
$ws.Range("D14").Value = "'" + '0.07112'
$ws.Range("E14").Value = '  +0.99%  '
$ws.Range("D15").Value = "'" + '6.511'
$ws.Range("E15").Value = '  +0.04%  '
$ws.Range("D16").Value = "'" + '91.39'
$ws.Range("E16").Value = '  -0.96%  '
$ws.Range("D17").Value = "'" + '1.003'
$ws.Range("E17").Value = '  +0.14%  '
$ws.Range("D18").Value = "'" + '0.000008695'
$ws.Range("E18").Value = '  -0.12%  '
$ws.Range("E19").Value = '  +0.15%  '
$ws.Range("E20").Value = '  -0.56%  '
$ws.Range("D21").Value = "'" + '26.982.50'
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("D22").Value = "'" + '5.295'
$ws.Range("E22").Value = '  -0.58%  '
$ws.Range("D23").Value = "'" + '10.60'
$ws.Range("E23").Value = '  -0.40%  '
$ws.Range("D24").Value = "'" + '2.043.42'
$ws.Range("E24").Value = '  +1.31%  '
$ws.Range("E25").Value = '  -0.20%  '
$ws.Range("D26").Value = "'" + '151.03'
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("D27").Value = "'" + '18.45'
$ws.Range("E27").Value = '  +0.33%  '
$ws.Range("D28").Value = "'" + '2.141'
$ws.Range("D29").Value = "'" + '5.267'
$ws.Range("E29").Value = '  -1.49%  '
$ws.Range("D30").Value = "'" + '116.58'
$ws.Range("E30").Value = '  +0.50%  '
$ws.Range("D31").Value = "'" + '0.08898'
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("D32").Value = "'" + '0.7605'
$ws.Range("E32").Value = '  +0.39%  '
$ws.Range("E33").Value = '  +0.75%  '
$ws.Range("D34").Value = "'" + '4.490'
$ws.Range("E34").Value = '  +0.83%  '
$ws.Range("D35").Value = "'" + '2.902'
$ws.Range("E35").Value = '  -0.50%  '
$ws.Range("E36").Value = '  +0.16%  '
$ws.Range("D37").Value = "'" + '1.091'
$ws.Range("E37").Value = '  -1.07%  '
$ws.Range("D38").Value = "'" + '0.05291'
$ws.Range("E38").Value = '  +0.82%  '
$ws.Range("D39").Value = "'" + '0.01948'
$ws.Range("E39").Value = '  -1.09%  '
$ws.Range("D40").Value = "'" + '2.971'
$ws.Range("E40").Value = '  +1.33%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = "'" + '7.169'
$ws.Range("E41").Value = '  -0.52%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = "'" + '0.5288'
$ws.Range("E42").Value = '  -0.75%  '
$ws.Range("D43").Value = "'" + '2.330'
$ws.Range("E43").Value = '  -3.98%  '
$ws.Range("D44").Value = "'" + '0.1657'
$ws.Range("E44").Value = '  -0.42%  '
$ws.Range("D45").Value = "'" + '8.433'
$ws.Range("E45").Value = '  -0.89%  '
$ws.Range("D46").Value = "'" + '0.4860'
$ws.Range("E46").Value = '  -2.49%  '
$ws.Range("D47").Value = "'" + '10.48'
$ws.Range("E47").Value = '  +1.53%  '
$ws.Range("E48").Value = '  +0.17%  '
$ws.Range("D49").Value = "'" + '103.44'
$ws.Range("E49").Value = '  -0.51%  '
$ws.Range("E50").Value = '  -0.54%  '
$ws.Range("D51").Value = "'" + '0.06296'
$ws.Range("E51").Value = '  +0.00%  '
